# edit.ps1 - Applies the "verso-4a" template refresh:
#   1. Updates the datetimeFigureOut "Date Placeholder" field found on every
#      slide-layout (master) from 23/03/2025 to 02/04/2025.
#   2. Repositions / resizes a handful of shapes on slide 1 and shrinks the
#      "[conteudo]" placeholder's font from 10pt to 9pt.

function EmuToPt($emu) {
    # Shape.Left/Top/Width/Height are backed by a 32-bit float (points), so a
    # plain EMU/12700 division can be rounded down by one EMU once it round
    # trips through that float. A tiny nudge keeps the stored EMU exact.
    return ($emu / 12700) + 0.00003
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update the date field text on every custom (slide) layout.
# ---------------------------------------------------------------------
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shp = $layout.Shapes.Item($si)
        if ($shp.HasTextFrame) {
            $isDatePlaceholder = $false
            try {
                if ($shp.PlaceholderFormat.Type -eq 16) {
                    $isDatePlaceholder = $true
                }
            } catch {
            }
            if ($isDatePlaceholder) {
                if ($shp.TextFrame.TextRange.Text -eq "23/03/2025") {
                    $shp.TextFrame.TextRange.Text = "02/04/2025"
                }
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Adjust shapes on slide 1 (looked up by their stable shape names).
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)
$shapes = $slide.Shapes

# "Caixa de Texto 2" -> the "[conteudo]" textbox: move/resize + shrink font.
$conteudo = $shapes.Item("Caixa de Texto 2")
$conteudo.Left = EmuToPt 328687
$conteudo.Top = EmuToPt 1645044
$conteudo.Width = EmuToPt 9101348
$conteudo.Height = EmuToPt 2268000
$conteudo.TextFrame.TextRange.Font.Size = 9

# "Retângulo 6" -> signature box 2 placeholder block.
$ret6 = $shapes.Item("Retângulo 6")
$ret6.Top = EmuToPt 3931955

# "CaixaDeTexto 3" -> CONTRATANTE block.
$caixa3 = $shapes.Item("CaixaDeTexto 3")
$caixa3.Top = EmuToPt 5525196

# "Retângulo 4" -> signature box.
$ret4 = $shapes.Item("Retângulo 4")
$ret4.Top = EmuToPt 3931955

# "CaixaDeTexto 5" -> LOCAL DO TREINAMENTO block.
$caixa5 = $shapes.Item("CaixaDeTexto 5")
$caixa5.Top = EmuToPt 5998727

# "Retângulo 16" -> signature box.
$ret16 = $shapes.Item("Retângulo 16")
$ret16.Top = EmuToPt 3931955

# "Retângulo 10" -> signature box.
$ret10 = $shapes.Item("Retângulo 10")
$ret10.Top = EmuToPt 3931955
